$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New emmeans / SE / df / lower.CL / upper.CL values per row (2-61)
$data = @(
    @(2, 59.601, 0.576, 2007, 58.472, 60.731),
    @(3, 62.554, 0.281, 2007, 62.003, 63.104),
    @(4, 62.152, 0.552, 2007, 61.069, 63.235),
    @(5, 61.642, 0.858, 2007, 59.959, 63.325),
    @(6, 61.999, 0.226, 2007, 61.555, 62.443),
    @(7, 109.75, 1.011, 2007, 107.767, 111.733),
    @(8, 113.831, 0.493, 2007, 112.865, 114.798),
    @(9, 111.037, 0.97, 2007, 109.134, 112.939),
    @(10, 109.199, 1.507, 2007, 106.244, 112.155),
    @(11, 109.174, 0.398, 2007, 108.394, 109.954),
    @(12, 127.408, 1.295, 2007, 124.867, 129.948),
    @(13, 134.185, 0.631, 2007, 132.947, 135.423),
    @(14, 133.408, 1.243, 2007, 130.97, 135.845),
    @(15, 131.467, 1.931, 2007, 127.68, 135.254),
    @(16, 132.791, 0.509, 2007, 131.792, 133.79),
    @(17, 97.572, 1.409, 2007, 94.80800000000001, 100.335),
    @(18, 104.263, 0.6870000000000001, 2007, 102.916, 105.61),
    @(19, 100.008, 1.352, 2007, 97.357, 102.66),
    @(20, 98.601, 2.1, 2007, 94.482, 102.72),
    @(21, 100.296, 0.554, 2007, 99.20999999999999, 101.383),
    @(22, 16.446, 0.499, 2007, 15.467, 17.426),
    @(23, 19.551, 0.243, 2007, 19.073, 20.028),
    @(24, 17.401, 0.479, 2007, 16.461, 18.341),
    @(25, 18.984, 0.744, 2007, 17.524, 20.444),
    @(26, 17.039, 0.196, 2007, 16.654, 17.424),
    @(27, 18.373, 0.258, 2007, 17.867, 18.879),
    @(28, 18.292, 0.126, 2007, 18.045, 18.538),
    @(29, 18.959, 0.247, 2007, 18.474, 19.445),
    @(30, 18.227, 0.384, 2007, 17.473, 18.981),
    @(31, 19.93, 0.101, 2007, 19.731, 20.129),
    @(32, 43.81, 0.42, 2007, 42.986, 44.634),
    @(33, 43.491, 0.205, 2007, 43.089, 43.893),
    @(34, 45.52, 0.403, 2007, 44.729, 46.311),
    @(35, 45.11, 0.627, 2007, 43.881, 46.339),
    @(36, 45.395, 0.165, 2007, 45.071, 45.72),
    @(37, 111.833, 0.887, 2007, 110.093, 113.573),
    @(38, 118.845, 0.432, 2007, 117.997, 119.693),
    @(39, 116.333, 0.851, 2007, 114.664, 118.003),
    @(40, 115.726, 1.322, 2007, 113.132, 118.319),
    @(41, 116.373, 0.349, 2007, 115.689, 117.058),
    @(42, 71.01300000000001, 1.016, 2007, 69.021, 73.005),
    @(43, 80.45099999999999, 0.495, 2007, 79.48, 81.422),
    @(44, 73.83499999999999, 0.974, 2007, 71.92400000000001, 75.746),
    @(45, 73.996, 1.514, 2007, 71.027, 76.965),
    @(46, 73.086, 0.399, 2007, 72.30200000000001, 73.869),
    @(47, 151.819, 1.235, 2007, 149.396, 154.241),
    @(48, 158.626, 0.602, 2007, 157.445, 159.807),
    @(49, 154.316, 1.185, 2007, 151.991, 156.64),
    @(50, 153.133, 1.841, 2007, 149.522, 156.744),
    @(51, 154.532, 0.486, 2007, 153.579, 155.484),
    @(52, 281.723, 1.302, 2007, 279.17, 284.277),
    @(53, 289.374, 0.635, 2007, 288.129, 290.619),
    @(54, 281.905, 1.249, 2007, 279.455, 284.356),
    @(55, 282.612, 1.941, 2007, 278.806, 286.419),
    @(56, 282.394, 0.512, 2007, 281.39, 283.398),
    @(57, 147.574, 0.643, 2007, 146.313, 148.834),
    @(58, 149.106, 0.313, 2007, 148.492, 149.721),
    @(59, 147.395, 0.617, 2007, 146.186, 148.604),
    @(60, 147.712, 0.958, 2007, 145.833, 149.59),
    @(61, 146.407, 0.253, 2007, 145.911, 146.902)
)

foreach ($row in $data) {
    $r = $row[0]
    $arr = New-Object 'object[,]' 1,5
    $arr[0,0] = $row[1]
    $arr[0,1] = $row[2]
    $arr[0,2] = $row[3]
    $arr[0,3] = $row[4]
    $arr[0,4] = $row[5]
    $ws.Range("B$r`:F$r").Value = $arr
}
